$d = $word.ActiveDocument

# Fix typo: "following" -> "folowing"
$d.Content.Find.Execute("following adaptive", $true, $false, $false, $false, $false, $true, 1, $false, "folowing adaptive", 2)

# Remove the paragraph that holds the inline picture (the image has been dropped from this page)
$d.Paragraphs(3).Range.Delete()

# Fix "(i) | The aquatic plants can get sufficient nutrients to grow well." -> "(4) | The aquatic plants can get sufficient nutrients to grow wei."
$d.Content.Find.Execute("(i) | The aquatic plants can get sufficient nutrients to grow well.", $true, $false, $false, $false, $false, $true, 1, $false, "(4) | The aquatic plants can get sufficient nutrients to grow wei.", 2)
